$d = $word.ActiveDocument

# Find the paragraph that holds "Expert methodology validated at highest
# judicial level" (the last bullet of the KEY ACHIEVEMENTS AND IMPACT /
# Impact section) so we can insert the two new bullet paragraphs right
# after it, before the TECHNICAL SKILLS heading.
$total = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $total; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Expert methodology validated at highest judicial level*") {
        $anchorIndex = $i
        break
    }
}

$anchorPara = $d.Paragraphs($anchorIndex)
$rng = $anchorPara.Range
$rng.Collapse(0)

# --- First new paragraph: plain bullet text. ---
$rng.InsertParagraphAfter()
$p1 = $d.Paragraphs($anchorIndex + 1)
$p1.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# --- Second new paragraph: bullet with a bold/colored "178%" run. ---
$rng2 = $p1.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
$p2 = $d.Paragraphs($anchorIndex + 2)
$p2.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Bold + color the "178%" portion of the paragraph just written.
$boldRng = $p2.Range.Duplicate
$boldRng.Find.Execute("178%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldRng.Font.Bold = $true
$boldRng.Font.Color = 5258796
